# "Se agregan archivos de apoyo"
# The first checklist table (Tabla1, rows 1:6) together with its
# explanatory "example" comments is removed. The remaining checklist
# table (Tabla13) slides up so only a single blank row separates it
# from the top of the sheet (it now lives at A2:D7 instead of A10:D15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first table's rows (1-8, i.e. the 6 table rows plus the
# 2 blank rows that keep a single blank row before the second table).
# Deleting the rows removes the embedded ListObject (Tabla1) along
# with its data/comments, and shifts the second table (Tabla13) up
# from A10:D15 to A2:D7.
$ws.Rows("1:8").Delete()

# Restore the previously selected cell (shifted up along with the data).
$ws.Range("B17").Select()
